$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "AEDB.CEA"
$ws.Range("B2").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("C2").Value = "Macrophages_rank"
$ws.Range("D2").Value = 0.0714487617452087
$ws.Range("E2").Value = 0.032286644487318
$ws.Range("F2").Value = 1.07406311617178
$ws.Range("G2").Value = 1.00820037896599
$ws.Range("H2").Value = 1.14422847043936
$ws.Range("I2").Value = 2.21295098576359
$ws.Range("J2").Value = 0.0271282739699829
$ws.Range("K2").Value = 0.199642650122291
$ws.Range("L2").Value = 0.176221496594084
$ws.Range("M2").Value = 2423
$ws.Range("N2").Value = 1021
$ws.Range("O2").Value = 57.8621543541065

# Row 3
$ws.Range("A3").Value = "AEDB.CEA"
$ws.Range("B3").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("C3").Value = "SMC_rank"
$ws.Range("D3").Value = -0.073653066853816
$ws.Range("E3").Value = 0.0324488422366725
$ws.Range("F3").Value = 0.928993936728941
$ws.Range("G3").Value = 0.871749821978802
$ws.Range("H3").Value = 0.989997029790184
$ws.Range("I3").Value = -2.26982110229424
$ws.Range("J3").Value = 0.0234334189896108
$ws.Range("K3").Value = 0.199901279178963
$ws.Range("L3").Value = 0.176392806125458
$ws.Range("M3").Value = 2423
$ws.Range("N3").Value = 1017
$ws.Range("O3").Value = 58.027238959967

# Row 4
$ws.Range("A4").Value = "AEDB.CEA"
$ws.Range("B4").Value = "MCP1_pg_ug_2015_rank"
$ws.Range("C4").Value = "VesselDensity_rank"
$ws.Range("D4").Value = -0.0570783394030029
$ws.Range("E4").Value = 0.0316332950527268
$ws.Range("F4").Value = 0.94452007332949
$ws.Range("G4").Value = 0.887737136109742
$ws.Range("H4").Value = 1.00493505637468
$ws.Range("I4").Value = -1.80437539964977
$ws.Range("J4").Value = 0.0714983494109482
$ws.Range("K4").Value = 0.19999607245993
$ws.Range("L4").Value = 0.174860521107099
$ws.Range("M4").Value = 2423
$ws.Range("N4").Value = 953
$ws.Range("O4").Value = 60.668592653735

# Row 5
$ws.Range("A5").Value = "AEDB.CEA"
$ws.Range("B5").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("C5").Value = "Macrophages_rank"
$ws.Range("D5").Value = 0.0870152953076361
$ws.Range("E5").Value = 0.0308829201560253
$ws.Range("F5").Value = 1.0909133654462
$ws.Range("G5").Value = 1.02683861024077
$ws.Range("H5").Value = 1.15898638699425
$ws.Range("I5").Value = 2.8175863832831
$ws.Range("J5").Value = 0.00493480323451057
$ws.Range("K5").Value = 0.251597553987997
$ws.Range("L5").Value = 0.229696776052227
$ws.Range("M5").Value = 2423
$ws.Range("N5").Value = 1021
$ws.Range("O5").Value = 57.8621543541065

# Row 6
$ws.Range("A6").Value = "AEDB.CEA"
$ws.Range("B6").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("C6").Value = "SMC_rank"
$ws.Range("D6").Value = -0.145018671668191
$ws.Range("E6").Value = 0.0307987745268994
$ws.Range("F6").Value = 0.865006141852292
$ws.Range("G6").Value = 0.814334332237226
$ws.Range("H6").Value = 0.918830995847314
$ws.Range("I6").Value = -4.70858577640915
$ws.Range("J6").Value = 0.00000285093361592816
$ws.Range("K6").Value = 0.261683291659197
$ws.Range("L6").Value = 0.23999009556813
$ws.Range("M6").Value = 2423
$ws.Range("N6").Value = 1017
$ws.Range("O6").Value = 58.027238959967

# Row 7
$ws.Range("A7").Value = "AEDB.CEA"
$ws.Range("B7").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("C7").Value = "VesselDensity_rank"
$ws.Range("D7").Value = -0.0476932344938231
$ws.Range("E7").Value = 0.030287754555999
$ws.Range("F7").Value = 0.953426220498102
$ws.Range("G7").Value = 0.898474248761471
$ws.Range("H7").Value = 1.01173913352148
$ws.Range("I7").Value = -1.57467052916198
$ws.Range("J7").Value = 0.115675193595248
$ws.Range("K7").Value = 0.243832301412757
$ws.Range("L7").Value = 0.220074053028109
$ws.Range("M7").Value = 2423
$ws.Range("N7").Value = 953
$ws.Range("O7").Value = 60.668592653735

# Row 8
$ws.Range("A8").Value = "AEDB.CEA"
$ws.Range("B8").Value = "MCP1_rank"
$ws.Range("C8").Value = "Macrophages_rank"
$ws.Range("D8").Value = 0.0971796405511474
$ws.Range("E8").Value = 0.0418936498686737
$ws.Range("F8").Value = 1.1020583302058
$ws.Range("G8").Value = 1.01518218277759
$ws.Range("H8").Value = 1.19636906929647
$ws.Range("I8").Value = 2.31967472053119
$ws.Range("J8").Value = 0.0207806871586204
$ws.Range("K8").Value = 0.114927287761619
$ws.Range("L8").Value = 0.0777393586759728
$ws.Range("M8").Value = 2423
$ws.Range("N8").Value = 497
$ws.Range("O8").Value = 79.4882377218324

# Row 9
$ws.Range("A9").Value = "AEDB.CEA"
$ws.Range("B9").Value = "MCP1_rank"
$ws.Range("C9").Value = "SMC_rank"
$ws.Range("D9").Value = -0.253030666527429
$ws.Range("E9").Value = 0.0438266435506658
$ws.Range("F9").Value = 0.776444070615248
$ws.Range("G9").Value = 0.712531691958533
$ws.Range("H9").Value = 0.84608923588575
$ws.Range("I9").Value = -5.77344386947891
$ws.Range("J9").Value = 0.0000000140596782157026
$ws.Range("K9").Value = 0.164109031683249
$ws.Range("L9").Value = 0.128839370572837
$ws.Range("M9").Value = 2423
$ws.Range("N9").Value = 495
$ws.Range("O9").Value = 79.5707800247627

# Row 10
$ws.Range("A10").Value = "AEDB.CEA"
$ws.Range("B10").Value = "MCP1_rank"
$ws.Range("C10").Value = "VesselDensity_rank"
$ws.Range("D10").Value = -0.042854069861196
$ws.Range("E10").Value = 0.0551328272180718
$ws.Range("F10").Value = 0.958051188409293
$ws.Range("G10").Value = 0.859921295133979
$ws.Range("H10").Value = 1.06737917156645
$ws.Range("I10").Value = -0.777287725363537
$ws.Range("J10").Value = 0.437383439023689
$ws.Range("K10").Value = 0.109490143209984
$ws.Range("L10").Value = 0.0712708360516142
$ws.Range("M10").Value = 2423
$ws.Range("N10").Value = 487
$ws.Range("O10").Value = 79.9009492364837
